$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 472.125
$ws.Range("I33").Value = 491.7143
$ws.Range("J33").Value = 335
$ws.Range("K33").Value = 491.7143
$ws.Range("L33").Value = 335
$ws.Range("M33").Value = -262.7143
$ws.Range("N33").Value = -793
# Row 40
$ws.Range("H40").Value = 40332.668
$ws.Range("I40").Value = 38499
$ws.Range("K40").Value = 38499
$ws.Range("M40").Value = -38324
# Row 70
$ws.Range("H70").Value = 1899.5454
$ws.Range("I70").Value = 1569
$ws.Range("J70").Value = 2296.2
$ws.Range("K70").Value = 4707
$ws.Range("L70").Value = 6888.599999999999
$ws.Range("M70").Value = -4437
$ws.Range("N70").Value = -7428.599999999999
# Row 73
$ws.Range("H73").Value = 1899.5454
$ws.Range("I73").Value = 1569
$ws.Range("J73").Value = 2296.2
$ws.Range("K73").Value = 4707
$ws.Range("L73").Value = 6888.599999999999
$ws.Range("M73").Value = -3771
$ws.Range("N73").Value = -8760.599999999999
# Row 76
$ws.Range("H76").Value = 1999.5
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 1999.5
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 86
$ws.Range("H86").Value = 200003470
$ws.Range("I86").Value = 250003200
$ws.Range("J86").Value = 4499.6665
$ws.Range("K86").Value = 250003200
$ws.Range("L86").Value = 4499.6665
$ws.Range("M86").Value = -250002077
$ws.Range("N86").Value = -6745.6665
# Row 89
$ws.Range("H89").Value = 200003470
$ws.Range("I89").Value = 250003200
$ws.Range("J89").Value = 4499.6665
$ws.Range("K89").Value = 1250016000
$ws.Range("L89").Value = 22498.3325
$ws.Range("M89").Value = -1250010384
$ws.Range("N89").Value = -33730.3325
# Row 92
$ws.Range("H92").Value = 57972068
$ws.Range("I92").Value = 70176520
$ws.Range("K92").Value = 70176520
$ws.Range("M92").Value = -70175272
# Row 107
$ws.Range("H107").Value = 40002040
$ws.Range("I107").Value = 62501216
$ws.Range("J107").Value = 3500.5557
$ws.Range("K107").Value = 62501216
$ws.Range("L107").Value = 3500.5557
$ws.Range("M107").Value = -62499296
$ws.Range("N107").Value = -7340.5557
# Row 116
$ws.Range("H116").Value = 5198.2
$ws.Range("J116").Value = 4518.4287
$ws.Range("L116").Value = 4518.4287
$ws.Range("N116").Value = -11402.4287

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 8891.421
$ws.Range("I45").Value = 12487.667
$ws.Range("J45").Value = 2726.4285
$ws.Range("K45").Value = 12487.667
$ws.Range("L45").Value = 2726.4285
$ws.Range("M45").Value = -12110.667
$ws.Range("N45").Value = -3480.4285
# Row 122
$ws.Range("H122").Value = 6862.411
$ws.Range("I122").Value = 4679.3125
$ws.Range("K122").Value = 14037.9375
$ws.Range("M122").Value = -11587.9375
# Row 132
$ws.Range("H132").Value = 5534.4526
$ws.Range("I132").Value = 4852.0835
$ws.Range("K132").Value = 14556.2505
$ws.Range("M132").Value = -12026.2505

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2999.9167
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 3999.875
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 3999.875
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -6245.875
# Row 89
$ws.Range("H89").Value = 2999.9167
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 3999.875
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 19999.375
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -31231.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5123.607
$ws.Range("I58").Value = 5506.2173
$ws.Range("J58").Value = 3363.6
$ws.Range("K58").Value = 5506.2173
$ws.Range("L58").Value = 3363.6
$ws.Range("M58").Value = -5303.2173
$ws.Range("N58").Value = -3769.6
# Row 74
$ws.Range("H74").Value = 64664
$ws.Range("J74").Value = 64664
$ws.Range("L74").Value = 64664
$ws.Range("N74").Value = -66412
# Row 77
$ws.Range("H77").Value = 64664
$ws.Range("J77").Value = 64664
$ws.Range("L77").Value = 193992
$ws.Range("N77").Value = -202728
# Row 134
$ws.Range("H134").Value = 5256.3716
$ws.Range("I134").Value = 4076.8147
$ws.Range("K134").Value = 12230.4441
$ws.Range("M134").Value = -9695.444100000001
# Row 136
$ws.Range("H136").Value = 5123.607
$ws.Range("I136").Value = 5506.2173
$ws.Range("J136").Value = 3363.6
$ws.Range("K136").Value = 16518.6519
$ws.Range("L136").Value = 10090.8
$ws.Range("M136").Value = -13968.6519
$ws.Range("N136").Value = -15190.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 103.65
$ws.Range("I2").Value = 44.75
$ws.Range("J2").Value = 142.91667
$ws.Range("K2").Value = 268.5
$ws.Range("L2").Value = 857.5000200000001
$ws.Range("M2").Value = -155.5
$ws.Range("N2").Value = -1083.50002
# Row 107
$ws.Range("H107").Value = 2221.7778
$ws.Range("J107").Value = 2114
$ws.Range("L107").Value = 6342
$ws.Range("N107").Value = -10182
# Row 131
$ws.Range("H131").Value = 12858910
$ws.Range("I131").Value = 12032880
$ws.Range("J131").Value = 13891448
$ws.Range("K131").Value = 36098640
$ws.Range("L131").Value = 41674344
$ws.Range("M131").Value = -36093600
$ws.Range("N131").Value = -41684424

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3681.7144
$ws.Range("I113").Value = 1842.2858
$ws.Range("J113").Value = 5521.143
$ws.Range("K113").Value = 1842.2858
$ws.Range("L113").Value = 5521.143
$ws.Range("M113").Value = 327.7141999999999
$ws.Range("N113").Value = -9861.143
# Row 132
$ws.Range("H132").Value = 7466.3335
$ws.Range("I132").Value = 7399
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 22197
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -19667
$ws.Range("N132").Value = -27560

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1459.8
$ws.Range("J93").Value = 1466.6666
$ws.Range("L93").Value = 1466.6666
$ws.Range("N93").Value = -3962.6666
# Row 136
$ws.Range("H136").Value = 3835170.8
$ws.Range("J136").Value = 9498.556
$ws.Range("L136").Value = 28495.668
$ws.Range("N136").Value = -33595.66800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1477.8572
$ws.Range("I113").Value = 1474.2222
$ws.Range("K113").Value = 4422.6666
$ws.Range("M113").Value = -2252.6666
# Row 126
$ws.Range("H126").Value = 4593.975
$ws.Range("J126").Value = 6266.875
$ws.Range("L126").Value = 18800.625
$ws.Range("N126").Value = -23740.625
# Row 132
$ws.Range("H132").Value = 3463.261
$ws.Range("I132").Value = 3502.6667
$ws.Range("K132").Value = 10508.0001
$ws.Range("M132").Value = -7978.000100000001
